# Adds the two new log rows (10 and 11) to the bottom of the tracking
# sheet, for the "receive" and "forward" actions that were logged.
#
# Columns: id | name | recentDate | recentTime | recentPlace |
#          previousDate | previousTime | previousPlace

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C and F hold dates formatted as plain text (e.g. "2025-01-08" or
# the "0001-01-01" placeholder for "no previous date"). Force those two
# columns to Text format *before* writing so Excel does not silently
# reinterpret the strings as real date serial numbers.
$ws.Range("C10:C11").NumberFormat = "@"
$ws.Range("F10:F11").NumberFormat = "@"

# Row 10
$ws.Cells.Item(10, 1).Value = "DOC-1735970018472"
$ws.Cells.Item(10, 2).Value = "Manvir"
$ws.Cells.Item(10, 3).Value = "2025-01-08"
$ws.Cells.Item(10, 4).Value = "14:02"
$ws.Cells.Item(10, 5).Value = "manvir"
$ws.Cells.Item(10, 6).Value = "0001-01-01"
$ws.Cells.Item(10, 7).Value = "13:30"
$ws.Cells.Item(10, 8).Value = "kdsfs"

# Row 11
$ws.Cells.Item(11, 1).Value = "DOC-1735970506282"
$ws.Cells.Item(11, 2).Value = "maa"
$ws.Cells.Item(11, 3).Value = "2025-01-14"
$ws.Cells.Item(11, 4).Value = "12:04"
$ws.Cells.Item(11, 5).Value = "Manvir"
$ws.Cells.Item(11, 6).Value = "2025-01-22"
$ws.Cells.Item(11, 7).Value = "12:02"
$ws.Cells.Item(11, 8).Value = "feds"
